$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value = "67e72c4acde4e894925703ff"
$ws.Range("V3").Value = "67e72c4acde4e894925703fd"
$ws.Range("V4").Value = "67e72c4acde4e89492570401"
$ws.Range("V5").Value = "67e72c49cde4e8949257024b"
$ws.Range("V6").Value = "67e72c49cde4e8949257024f"
$ws.Range("V7").Value = "67e72c49cde4e8949257024d"
$ws.Range("V8").Value = "67e72c49cde4e89492570251"
$ws.Range("V9").Value = "67e72c4acde4e89492570515"
$ws.Range("V10").Value = "67e72c4acde4e89492570517"
$ws.Range("V11").Value = "67e72c4acde4e89492570519"
$ws.Range("V12").Value = "67e72c4acde4e89492570667"
$ws.Range("V13").Value = "67e72c4acde4e89492570665"
$ws.Range("V14").Value = "67e72c49cde4e894925703d0"
$ws.Range("V15").Value = "67e72c49cde4e894925703d4"
$ws.Range("V16").Value = "67e72c49cde4e894925703ce"
$ws.Range("V17").Value = "67e72c49cde4e894925703d2"
$ws.Range("V18").Value = "67e72c4acde4e8949257056f"
$ws.Range("V19").Value = "67e72c4acde4e89492570571"
$ws.Range("V20").Value = "67e72c49cde4e894925702cd"
$ws.Range("V21").Value = "67e72c49cde4e894925702cb"
$ws.Range("V22").Value = "67e72c49cde4e894925702c9"
$ws.Range("V23").Value = "67e72c48cde4e894925701fb"
$ws.Range("V24").Value = "67e72c48cde4e894925701fd"
$ws.Range("V25").Value = "67e72c48cde4e894925701ff"
$ws.Range("V26").Value = "67e72c49cde4e894925703a1"
$ws.Range("V27").Value = "67e72c49cde4e8949257039f"
$ws.Range("V28").Value = "67e72c49cde4e894925703a5"
$ws.Range("V29").Value = "67e72c4acde4e8949257045b"
$ws.Range("V30").Value = "67e72c4acde4e8949257045f"
$ws.Range("V31").Value = "67e72c4acde4e89492570465"
$ws.Range("V32").Value = "67e72c4acde4e8949257045d"
$ws.Range("V33").Value = "67e72c4acde4e89492570461"
$ws.Range("V34").Value = "67e72c49cde4e8949257028b"
$ws.Range("V35").Value = "67e72c49cde4e8949257028d"
$ws.Range("V36").Value = "67e72c49cde4e8949257028d"
$ws.Range("V37").Value = "67e72c4acde4e89492570428"
$ws.Range("V38").Value = "67e72c4acde4e89492570424"
$ws.Range("V39").Value = "67e72c4acde4e89492570420"
$ws.Range("V40").Value = "67e72c4acde4e89492570422"
$ws.Range("V41").Value = "67e72c4acde4e89492570426"
$ws.Range("V42").Value = "67e72c4acde4e8949257053a"
$ws.Range("V43").Value = "67e72c4acde4e89492570536"
$ws.Range("V44").Value = "67e72c4acde4e89492570538"
$ws.Range("V45").Value = "67e72c4acde4e89492570534"
$ws.Range("V46").Value = "67e72c4acde4e894925704a3"
$ws.Range("V47").Value = "67e72c4acde4e894925704a5"
$ws.Range("V48").Value = "67e72c4acde4e894925704a9"
$ws.Range("V49").Value = "67e72c4acde4e894925704a7"
$ws.Range("V50").Value = "67e72c49cde4e89492570276"
$ws.Range("V51").Value = "67e72c49cde4e89492570278"
$ws.Range("V52").Value = "67e72c4acde4e8949257067a"
$ws.Range("V53").Value = "67e72c4acde4e8949257067c"
$ws.Range("V54").Value = "67e72c4acde4e8949257067e"
$ws.Range("V55").Value = "67e72c4acde4e894925704e9"
$ws.Range("V56").Value = "67e72c4acde4e894925704e5"
$ws.Range("V57").Value = "67e72c4acde4e894925704e3"
$ws.Range("V58").Value = "67e72c4acde4e894925705c5"
$ws.Range("V59").Value = "67e72c4acde4e894925705c3"
$ws.Range("V60").Value = "67e72c4acde4e894925705c7"
$ws.Range("V61").Value = "67e72c4acde4e894925705c9"
$ws.Range("V62").Value = "67e72c49cde4e89492570317"
$ws.Range("V63").Value = "67e72c49cde4e89492570319"
$ws.Range("V64").Value = "67e72c49cde4e89492570315"
$ws.Range("V65").Value = "67e72c49cde4e89492570351"
$ws.Range("V66").Value = "67e72c49cde4e89492570355"
$ws.Range("V67").Value = "67e72c49cde4e8949257033a"
$ws.Range("V68").Value = "67e72c49cde4e89492570338"
$ws.Range("V69").Value = "67e72c4acde4e894925705fe"
$ws.Range("V70").Value = "67e72c4acde4e894925705fa"
$ws.Range("V71").Value = "67e72c4acde4e894925705f6"
$ws.Range("V72").Value = "67e72c49cde4e89492570372"
$ws.Range("V73").Value = "67e72c49cde4e8949257036e"
$ws.Range("V74").Value = "67e72c49cde4e89492570374"
$ws.Range("V75").Value = "67e72c49cde4e89492570370"
$ws.Range("V76").Value = "67e72c4bcde4e894925706ee"
$ws.Range("V77").Value = "67e72c4bcde4e894925706f2"
$ws.Range("V78").Value = "67e72c4bcde4e894925706f0"
$ws.Range("V79").Value = "67e72c4bcde4e894925706f4"
$ws.Range("V80").Value = "67e72c49cde4e894925702fd"
$ws.Range("V81").Value = "67e72c4bcde4e8949257071b"
$ws.Range("V82").Value = "67e72c4bcde4e89492570719"
$ws.Range("V83").Value = "67e72c4bcde4e8949257071d"
$ws.Range("V84").Value = "67e72c4acde4e89492570588"
$ws.Range("V85").Value = "67e72c4acde4e8949257058a"
$ws.Range("V86").Value = "67e72c4acde4e894925705a1"
$ws.Range("V87").Value = "67e72c4acde4e89492570699"
$ws.Range("V88").Value = "67e72c4acde4e894925705b8"
$ws.Range("V89").Value = "67e72c49cde4e894925702b0"
$ws.Range("V90").Value = "67e72c49cde4e894925702ae"
$ws.Range("V91").Value = "67e72c49cde4e89492570222"
$ws.Range("V92").Value = "67e72c49cde4e89492570224"
$ws.Range("V93").Value = "67e72c49cde4e8949257021e"
$ws.Range("V94").Value = "67e72c49cde4e89492570220"
$ws.Range("V95").Value = "67e72c4acde4e89492570636"
$ws.Range("V96").Value = "67e72c4acde4e89492570638"
$ws.Range("V97").Value = "67e72c4acde4e8949257063c"
$ws.Range("V98").Value = "67e72c4acde4e8949257063a"
$ws.Range("V99").Value = "67e72c4acde4e89492570623"
$ws.Range("V100").Value = "67e72c4acde4e89492570621"
$ws.Range("V101").Value = "67e72c4acde4e89492570504"
$ws.Range("V102").Value = "67e72c4bcde4e894925706b9"
$ws.Range("V103").Value = "67e72c4bcde4e894925706c1"
$ws.Range("V104").Value = "67e72c4bcde4e894925706c1"
$ws.Range("V105").Value = "67e72c4bcde4e894925706bb"
$ws.Range("V106").Value = "67e72c4bcde4e894925706bd"
$ws.Range("V107").Value = "67e72c4bcde4e894925706a8"
$ws.Range("V108").Value = "67e72c4bcde4e89492570746"
$ws.Range("V109").Value = "67e72c4bcde4e89492570744"
$ws.Range("V110").Value = "67e72c4bcde4e8949257078d"
$ws.Range("V111").Value = "67e72c4bcde4e89492570791"
$ws.Range("V113").Value = "67e72c4bcde4e8949257086f"
$ws.Range("V115").Value = "67e72c4bcde4e894925707a6"
$ws.Range("V116").Value = "67e72c4bcde4e89492570806"
$ws.Range("V121").Value = "67e72c4bcde4e89492570823"
$ws.Range("V122").Value = "67e72c4bcde4e89492570811"
$ws.Range("V127").Value = "67e72c4bcde4e89492570833"
$ws.Range("V128").Value = "67e72c4bcde4e894925707c8"
$ws.Range("V129").Value = "67e72c4bcde4e89492570845"
$ws.Range("V130").Value = "67e72c4bcde4e89492570852"
$ws.Range("V131").Value = "67e72c4bcde4e89492570778"
$ws.Range("V132").Value = "67e72c4bcde4e89492570862"
$ws.Range("V133").Value = "67e72c4bcde4e894925707ee"
